$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A14").Value = "13_Case2_"
$ws.Range("B14").Value = "[('A', 'Spades'), ('A', 'Hearts'), ('A', 'Clubs'), ('2', 'Spades'), ('4', 'Hearts'), ('5', 'Clubs'), ('6', 'Diamonds'), ('7', 'Spades'), ('8', 'Hearts'), ('9', 'Clubs'), ('9', 'Diamonds'), ('10','Clubs'), ('10', 'Diamonds')]"
$ws.Range("C14").Value = "Needed to improve test_dim function with tuple, so each straight in Case 2 is registered properly."

$ws.Range("C14").Select()
